# PSP_Sheet_2조.xlsx — fill in the Oct 24 2019 time-log entry on the
# "작성자명" sheet (row 17) and move the active selection to F21,
# matching the author's re-upload of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Row 17: Date / Start / Stop / Interruption / Delta / Activity ---
$ws.Range("A17").Value = 43762                     # 2019-10-24
$ws.Range("B17").Value = 0.375                     # 09:00
$ws.Range("C17").Value = 0.45833333333333331       # 11:00
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 120

$activity = "Initial Data, Use Case Spec, DB 설계 등 각자 작업물 Review"
$ws.Range("F17").Value = $activity

# Second half of the activity text ( 등 각자 작업물 Review) keeps the
# 돋움/10pt run formatting used throughout this column's rich-text cells.
$splitAt = ("Initial Data, Use Case Spec, DB 설계").Length
$run = $ws.Range("F17").Characters($splitAt + 1, $activity.Length - $splitAt)
$run.Font.Name = "돋움"
$run.Font.Size = 10

# --- Move the active selection to F21 ---
[void]$ws.Range("F21").Select()
